# Auto-generated Excel COM-interop script updating the cryptos price table
# Applies the GitHub Actions "Updated cryptos list" data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.ClearFormats()
}

Set-TextValue "D2" "28.657.28"
$ws.Range("E2").Value = "  +1.46%  "

Set-TextValue "D3" "1.866.45"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  +0.45%  "

Set-TextValue "D5" "326.68"
$ws.Range("E5").Value = "  -1.17%  "

Set-TextValue "D6" "1.006"
$ws.Range("E6").Value = "  +0.43%  "

Set-TextValue "D7" "0.4629"
$ws.Range("E7").Value = "  +0.65%  "

Set-TextValue "D8" "0.3911"
$ws.Range("E8").Value = "  +1.38%  "

Set-TextValue "D9" "0.07895"
$ws.Range("E9").Value = "  +0.62%  "

Set-TextValue "D10" "0.9699"
$ws.Range("E10").Value = "  +0.63%  "

Set-TextValue "D11" "22.29"
$ws.Range("E11").Value = "  +2.14%  "

Set-TextValue "D12" "1.846.76"
$ws.Range("E12").Value = "  +1.03%  "

Set-TextValue "D13" "5.727"
$ws.Range("E13").Value = "  +0.08%  "

Set-TextValue "D14" "6.931"
$ws.Range("E14").Value = "  +0.29%  "

Set-TextValue "D15" "0.06931"
$ws.Range("E15").Value = "  +1.10%  "

Set-TextValue "D16" "88.54"
$ws.Range("E16").Value = "  +2.25%  "

Set-TextValue "D17" "1.007"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("E18").Value = "  +1.50%  "

$ws.Range("E19").Value = "  +0.44%  "

$ws.Range("E20").Value = "  +0.32%  "

Set-TextValue "D21" "28.683.24"
$ws.Range("E21").Value = "  +1.45%  "

Set-TextValue "D22" "5.320"
$ws.Range("E22").Value = "  -0.20%  "

Set-TextValue "D23" "11.07"
$ws.Range("E23").Value = "  +0.79%  "

Set-TextValue "D24" "2.127"
$ws.Range("E24").Value = "  -1.31%  "

Set-TextValue "D25" "2.124.56"
$ws.Range("E25").Value = "  +3.45%  "

Set-TextValue "D26" "155.22"
$ws.Range("E26").Value = "  +1.29%  "

Set-TextValue "D27" "19.31"
$ws.Range("E27").Value = "  +0.69%  "

Set-TextValue "D28" "5.771"
$ws.Range("E28").Value = "  -0.80%  "

Set-TextValue "D29" "1.993"
$ws.Range("E29").Value = "  +1.12%  "

Set-TextValue "D30" "119.18"
$ws.Range("E30").Value = "  +2.14%  "

$ws.Range("E31").Value = "  +0.32%  "

Set-TextValue "D32" "0.9373"
$ws.Range("E32").Value = "  -0.39%  "

Set-TextValue "D33" "5.318"
$ws.Range("E33").Value = "  +0.63%  "

Set-TextValue "D35" "3.348"
$ws.Range("E35").Value = "  -2.75%  "

Set-TextValue "D36" "0.05817"
$ws.Range("E36").Value = "  -3.59%  "

Set-TextValue "D37" "0.02114"
$ws.Range("E37").Value = "  -2.10%  "

Set-TextValue "D38" "1.157"
$ws.Range("E38").Value = "  +0.28%  "

Set-TextValue "D39" "7.894"
$ws.Range("E39").Value = "  +4.51%  "

Set-TextValue "D40" "0.5654"
$ws.Range("E40").Value = "  +0.67%  "

Set-TextValue "D41" "9.931"
$ws.Range("E41").Value = "  -0.66%  "

Set-TextValue "D42" "0.1776"
$ws.Range("E42").Value = "  -0.32%  "

Set-TextValue "D45" "11.69"
$ws.Range("E45").Value = "  -0.14%  "

Set-TextValue "D46" "0.5316"
$ws.Range("E46").Value = "  +0.54%  "

Set-TextValue "D47" "1.140"
$ws.Range("E47").Value = "  -8.10%  "

Set-TextValue "D48" "1.847"
$ws.Range("E48").Value = "  +0.50%  "

Set-TextValue "D49" "113.55"
$ws.Range("E49").Value = "  +0.50%  "

$ws.Range("E50").Value = "  +1.38%  "

Set-TextValue "D51" "1.006"
$ws.Range("E51").Value = "  +0.47%  "

# Row 43/44: RenderToken and Cronos swapped position in the ranking
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D43" "0.07252"
$ws.Range("E43").Value = "  +3.26%  "

$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "2.201"
$ws.Range("E44").Value = "  -2.36%  "
